# Update the results of Case 5-122 (380 kV) in the active sheet.
# Only columns B, C, D, E, F, I, J, O for rows 2..25 change; all other
# cells (A, G, H, K, L, M, N and the header row) stay exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(1.287629211722333,1.142681825063619,1.053468655672873,1.017061812461407,1.011013425886119,1.05297786675203,1.237697270622903,1.598144478682116,1.861792553738155,1.981462728161262,2.026738800614794,2.01698962365856,1.985188441541936,1.96570396486868,1.853966297756358,1.785349379317609,1.745857971652526,1.732482677738119,1.79265634837509,1.99453033546996,2.126229946432659,2.055961819632387,1.789353000266146,1.500833966780817)
$Cvals = @(0.3402632080236003,0.2985870368923429,0.2729121377406045,0.2624285278010063,0.2606864900856465,0.2727708358056873,0.3259113889523633,0.4294182518974594,0.5050141610020091,0.5393025812848578,0.5522717415465195,0.5494792837983482,0.5403698687111955,0.5347881002750796,0.5027712521620629,0.4831037285040907,0.4717820651450779,0.4679471462199558,0.4851983501435484,0.5430459403461896,0.580764249202673,0.5606415869202124,0.4842514169519632,0.4014944072816888)
$Dvals = @(0.235289304122503,0.2286407175799638,0.2246253972989507,0.223006075780205,0.222738216264375,0.2246034897530933,0.2329830487546332,0.2499421218842457,0.2627187381350211,0.2685991201601752,0.270835580772669,0.270353490081348,0.2687829216778255,0.2678221611054568,0.2623358037374146,0.2589874889381605,0.2570680568401826,0.2564192782829195,0.2593432581891761,0.2692439737597141,0.2757710730730878,0.2722823167964208,0.2591823975887166,0.2452983013398011)
$Evals = @(0.1776475392211836,0.1732199729131949,0.1705833462974233,0.1695295362741263,0.1693557994635775,0.1705690506382957,0.1761039334156393,0.1876068884387436,0.1964537076048174,0.2005643137564022,0.2021332615412703,0.2017948118587043,0.2006931447238145,0.2000199490965286,0.1961868000946083,0.1938573244432575,0.1925255821061143,0.1920760711396881,0.1941044619716052,0.201016396171255,0.2056057156967128,0.2031497362579699,0.1939927077186709,0.1844255579438823)
$Fvals = @(1.091381633918594,1.09542747142136,1.098777232410193,1.100359546011816,1.100635398054251,1.098797692910395,1.09259673831346,1.08732428579772,1.087678212699885,1.08876349448434,1.08930784876388,1.089184672650958,1.088805602975881,1.088590796051164,1.087625926162858,1.08727107552923,1.08715394419643,1.087129208231175,1.087299845560025,1.08891332123757,1.090745489015333,1.089696314427471,1.087286568044988,1.088010284904719)
$Ivals = @(0.4404514370260948,0.4510568956887351,0.4580034490540932,0.4609432486398219,0.4614379759974234,0.4580426551914947,0.4440177853716083,0.4199777512684548,0.4044438692803851,0.3978429548354985,0.3954106240007942,0.3959314715530464,0.3976414948729019,0.3986977073523406,0.4048846520719387,0.4087996238564724,0.4110952086980362,0.4118799657276639,0.4083783331028048,0.3971373904914532,0.3901831775703677,0.3938587522606678,0.4085686590790054,0.4261087158098515)
$Jvals = @(0.1793734731068,0.1753624964209948,0.1730093420532199,0.1720779490768649,0.1719249544422041,0.1729966694979836,0.1779677245160158,0.188587416441905,0.1969246061329102,0.2008344094957266,0.2023318455441654,0.2020085947674062,0.2009572660270891,0.2003154957382947,0.1966714515952503,0.1944659801295643,0.1932084785111954,0.192784603129013,0.1946996148461864,0.2012656084578737,0.205655258396277,0.2033034070149426,0.1945939560335148,0.1856208368358523)
$Ovals = @(2.400617408154631,2.424267857758139,2.440971406002291,2.448325514310213,2.449579670792332,2.441068372074341,2.40831833560614,2.361474758672671,2.337745632850044,2.329291069413131,2.326427615556128,2.327029248373123,2.329048705895275,2.330329761496529,2.338345383004167,2.3438631787308,2.347257017635542,2.348443874343275,2.343253003588728,2.328446353231413,2.320740843004842,2.324672490841067,2.343528173659621,2.372277102683881)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $Bvals[$i]
    $ws.Cells.Item($row, 3).Value = $Cvals[$i]
    $ws.Cells.Item($row, 4).Value = $Dvals[$i]
    $ws.Cells.Item($row, 5).Value = $Evals[$i]
    $ws.Cells.Item($row, 6).Value = $Fvals[$i]
    $ws.Cells.Item($row, 9).Value = $Ivals[$i]
    $ws.Cells.Item($row, 10).Value = $Jvals[$i]
    $ws.Cells.Item($row, 15).Value = $Ovals[$i]
}
